$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header format (bold font + grey fill, style index 2)
# onto the new legend-table header cells J1:K1, then fill in their text.
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "Progress"
$ws.Range("K1").Value = "meaning"

# Copy the percentage number format (style index 3, used in column E)
# onto J2:J6 before filling in the legend values.
$ws.Range("E2").Copy()
$ws.Range("J2:J6").PasteSpecial(-4122)  # xlPasteFormats

# Legend rows mapping each Progress percentage to its meaning
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "Not Started"

$ws.Range("J3").Value = 0.6
$ws.Range("K3").Value = "Major review needed"

$ws.Range("J4").Value = 0.85
$ws.Range("K4").Value = "Minor questions"

$ws.Range("J5").Value = 0.9
$ws.Range("K5").Value = "Not noted"

$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "Done"

$excel.CutCopyMode = 0

# Update the selected cell to match the saved workbook state
$ws.Range("K3").Select()
